# The dataset rows (2-21) are reordered: each destination row receives the
# D/J/K/L/M/O/P values that used to live in a different source row. Columns
# A, B, C, E, F, G, H, I, N, Q, R are identical across all rows so they do
# not need to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values captured from the sheet
# BEFORE any writes, since this is a permutation of existing row data).
$mapping = @{
    2  = 4
    3  = 19
    4  = 13
    5  = 7
    6  = 15
    7  = 11
    8  = 14
    9  = 18
    10 = 12
    11 = 16
    12 = 21
    13 = 3
    14 = 20
    15 = 5
    16 = 2
    17 = 8
    18 = 9
    19 = 10
    20 = 17
    21 = 6
}

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot the original values for the columns that move, for every row.
$original = @{}
foreach ($r in 2..21) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $original[$r] = $rowVals
}

# Write the permuted values into each destination row.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcVals[$col]
    }
}
